$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (font/border/alignment) template from the last existing row (845)
$fmtSrc = $ws.Range("A845:G845")

# Row 846
$ws.Cells.Item(846,1).Value = "MH"
$ws.Cells.Item(846,2).Value = "RGH1279_PNE_P40"
$ws.Cells.Item(846,3).Value = "11-Dec-2025 11:41 PM"
$ws.Cells.Item(846,4).Value = "FAIL"
$ws.Cells.Item(846,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(846,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(846,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrc.Copy()
$ws.Range("A846:G846").PasteSpecial(-4122)
$ws.Rows.Item(846).RowHeight = 104

# Row 847
$ws.Cells.Item(847,1).Value = "UW"
$ws.Cells.Item(847,2).Value = "AGA195_AGR_P40"
$ws.Cells.Item(847,3).Value = "11-Dec-2025 10:11 PM"
$ws.Cells.Item(847,4).Value = "FAIL"
$ws.Cells.Item(847,5).Value = "1. CSFB Call (pass/fail)`n2. Peak Rank - 5G"
$ws.Cells.Item(847,6).Value = "1. Static CSFB MO`n2. Static DL"
$ws.Cells.Item(847,7).Value = "1. CSFB MO – As per Bharti acceptance criteria, a minimum of 3 successful call setups are required without any blocked call. Kindly perform at least 3 successful MO  calls.`n2. Peak Rank is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area."
$fmtSrc.Copy()
$ws.Range("A847:G847").PasteSpecial(-4122)
$ws.Rows.Item(847).RowHeight = 104

# Row 848
$ws.Cells.Item(848,1).Value = "BH"
$ws.Cells.Item(848,2).Value = "BHPUY-02_PAT_P40"
$ws.Cells.Item(848,3).Value = "11-Dec-2025 9:16 PM"
$ws.Cells.Item(848,4).Value = "FAIL"
$ws.Cells.Item(848,5).Value = "1. SCG addition after VoLTE call released`n2. Peak PDCP DL Throughput`n3. Average PDCP DL Throughput`n4. Median PDCP DL​ Throughput`n5. Downlink Peak MCS - 5G`n6. Peak PDCP UL Throughput`n7. Avg PDCP UL Throughput`n8. Median PDCP UL Throughput`n9. UE Steering (Idle) : Non anchor/anchor to preferred anchor"
$ws.Cells.Item(848,6).Value = "1. Static VoLTE MO`n2. Static DL`n3. Static DL`n4. Mobility DL`n5. Static DL`n6. Static UL`n7. Static UL`n8. Mobility UL`n9. Static Idle"
$ws.Cells.Item(848,7).Value = "1. VoLTE Long Call MO – The SCG count after VoLTE call release should be equal to or greater than the total number of calls in the logfile. The KPI has failed because the SCG count is lower than the number of calls. Kindly exclude the existing logfile. While creating a new logfile, ensure the same site is serving in 4G and keep test files downloading in the background during static tests to support SCG addition.`n2. Peak PDCP DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP DL Throughput in the NR tab.`n3. Average PDCP DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP DL Throughput in the NR tab.`n4. The Median PDCP DL Throughput is reported as 0. Kindly add or exclude a logfile in the DL drive so the median value can update. It is recommended to add a new logfile and collect maximum throughput samples in a good coverage area.`n5. Peak MCS is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.`n6. Peak PDCP UL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP UL Throughput in the NR tab.`n7. Average PDCP UL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP UL Throughput in the NR tab.`n8. The Median PDCP UL Throughput is reported as 0. Kindly add or exclude a logfile in the DL drive so the median value can update. It is recommended to add a new logfile and collect maximum throughput samples in a good coverage area.`n9. For sites with NOKIA OEM, validate using Drive Idle, and for other OEMs, validate using Static Idle. In both Drive and Static Idle, the UE should latch from NR to LTE and from LTE to NR. In LTE, the UE should latch on the band that corresponds to the configured anchor layer."
$fmtSrc.Copy()
$ws.Range("A848:G848").PasteSpecial(-4122)
$ws.Rows.Item(848).RowHeight = 104

# Row 849
$ws.Cells.Item(849,1).Value = "OR"
$ws.Cells.Item(849,2).Value = "ORBHU-952_BHU_P41"
$ws.Cells.Item(849,3).Value = "11-Dec-2025 7:23 PM"
$ws.Cells.Item(849,4).Value = "FAIL"
$ws.Cells.Item(849,5).Value = "1. Video Streaming"
$ws.Cells.Item(849,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(849,7).Value = "1. While performing the YouTube test for both sectors, please ensure that the video is playing successfully in the script before saving the log file."
$fmtSrc.Copy()
$ws.Range("A849:G849").PasteSpecial(-4122)
$ws.Rows.Item(849).RowHeight = 104

# Row 850
$ws.Cells.Item(850,1).Value = "AP"
$ws.Cells.Item(850,2).Value = "HY9422_5g_HYD_P40"
$ws.Cells.Item(850,3).Value = "11-Dec-2025 6:29 PM"
$ws.Cells.Item(850,4).Value = "FAIL"
$ws.Cells.Item(850,5).Value = "1. Peak Rank - 5G`n2. Video Streaming  (ms)"
$ws.Cells.Item(850,6).Value = "1. Static DL`n2. Static Yotube Streaming"
$ws.Cells.Item(850,7).Value = "1. Peak Rank is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.`n2. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrc.Copy()
$ws.Range("A850:G850").PasteSpecial(-4122)
$ws.Rows.Item(850).RowHeight = 104

# Row 851
$ws.Cells.Item(851,1).Value = "MH"
$ws.Cells.Item(851,2).Value = "RTG2679_KLP_P40"
$ws.Cells.Item(851,3).Value = "11-Dec-2025 6:05 PM"
$ws.Cells.Item(851,4).Value = "FAIL"
$ws.Cells.Item(851,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(851,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(851,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrc.Copy()
$ws.Range("A851:G851").PasteSpecial(-4122)
$ws.Rows.Item(851).RowHeight = 104

# Row 852
$ws.Cells.Item(852,1).Value = "BH"
$ws.Cells.Item(852,2).Value = "BHCHA-73_PAT_P40"
$ws.Cells.Item(852,3).Value = "11-Dec-2025 6:15 PM"
$ws.Cells.Item(852,4).Value = "FAIL"
$ws.Cells.Item(852,5).Value = "1. Peak PDCP DL Throughput`n2. Average PDCP DL Throughput`n3. Median PDCP DL​ Throughput`n4. Peak PDCP UL Throughput`n5. Avg PDCP UL Throughput`n6. Peak PUSCH UL Throughput`n7. Median PDCP UL Throughput`n8. Uplink Peak MCS - 4G`n9. Peak Uplink PRB Allocation - 4G`n10. Ping/Round trip time(ms)"
$ws.Cells.Item(852,6).Value = "1. Static DL`n2. Static DL`n3. Mobility DL`n4. Static UL`n5. Static UL`n6. Static UL`n7. Mobility UL`n8. Static UL`n9. Static UL`n10. Static Ping"
$ws.Cells.Item(852,7).Value = "1. Peak PDCP DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP DL Throughput in the NR tab.`n2. Average PDCP DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP DL Throughput in the NR tab.`n3. The Median PDCP DL Throughput is reported as 0. Kindly add or exclude a logfile in the DL drive so the median value can update. It is recommended to add a new logfile and collect maximum throughput samples in a good coverage area.`n4. Peak PDCP UL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP UL Throughput in the NR tab.`n5. Average PDCP UL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the reporting of MR-DC DRB PDCP UL Throughput in the NR tab.`n6. Peak PUSCH DL Throughput is not meeting the acceptance criteria, and if the value is 0, it indicates that it was not recorded in the logfile. Kindly exclude the logfile and create a new one, and verify the maximum value of PUSCH Throughput in the NR tab.`n7. The Median PDCP UL Throughput is reported as 0. Kindly add or exclude a logfile in the DL drive so the median value can update. It is recommended to add a new logfile and collect maximum throughput samples in a good coverage area.`n8. Peak MCS is not meeting the acceptance criteria. Kindly redo the test and verify that the value meets the required threshold. To achieve the desired MCS, perform the test in the main lobe of the cell within a good coverage area.`n9. PRB is not reporting in your logfile. Kindly redo the test and verify the reporting of 4G PRB allocation.`n10. Ping is not meeting the acceptance criteria. The average ping value across all logfiles should be less than 50 ms. Kindly exclude the logfile where the average value exceeds 50 ms and redo the test.”"
$fmtSrc.Copy()
$ws.Range("A852:G852").PasteSpecial(-4122)
$ws.Rows.Item(852).RowHeight = 104

# Row 853
$ws.Cells.Item(853,1).Value = "MH"
$ws.Cells.Item(853,2).Value = "RTG3272_KLP_P40"
$ws.Cells.Item(853,3).Value = "11-Dec-2025 11:18 PM"
$ws.Cells.Item(853,4).Value = "FAIL"
$ws.Cells.Item(853,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(853,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(853,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrc.Copy()
$ws.Range("A853:G853").PasteSpecial(-4122)
$ws.Rows.Item(853).RowHeight = 104

# Row 854
$ws.Cells.Item(854,1).Value = "MH"
$ws.Cells.Item(854,2).Value = "BID6128_AUR_P40"
$ws.Cells.Item(854,3).Value = "11-Dec-2025 3:57 PM"
$ws.Cells.Item(854,4).Value = "FAIL"
$ws.Cells.Item(854,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(854,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(854,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrc.Copy()
$ws.Range("A854:G854").PasteSpecial(-4122)
$ws.Rows.Item(854).RowHeight = 104

# Row 855
$ws.Cells.Item(855,1).Value = "GJ"
$ws.Cells.Item(855,2).Value = "GU4999_GJ_P40"
$ws.Cells.Item(855,3).Value = "11-Dec-2025 3:32 PM"
$ws.Cells.Item(855,4).Value = "FAIL"
$ws.Cells.Item(855,5).Value = "1. Serving SSB beam steering"
$ws.Cells.Item(855,6).Value = "1. Mobility DL"
$ws.Cells.Item(855,7).Value = "1. Kindly add drive coverage in the failed sector and verify that the Beam Index servings are meeting the acceptance criteria."
$fmtSrc.Copy()
$ws.Range("A855:G855").PasteSpecial(-4122)
$ws.Rows.Item(855).RowHeight = 104

# Row 856
$ws.Cells.Item(856,1).Value = "BH"
$ws.Cells.Item(856,2).Value = "BHRIG-07_PAT_P40"
$ws.Cells.Item(856,3).Value = "11-Dec-2025 3:07 PM"
$ws.Cells.Item(856,4).Value = "FAIL"
$ws.Cells.Item(856,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(856,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(856,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrc.Copy()
$ws.Range("A856:G856").PasteSpecial(-4122)
$ws.Rows.Item(856).RowHeight = 104

# Row 857
$ws.Cells.Item(857,1).Value = "KK"
$ws.Cells.Item(857,2).Value = "NGGHI2_DAV_P40"
$ws.Cells.Item(857,3).Value = "11-Dec-2025 12:31 AM"
$ws.Cells.Item(857,4).Value = "FAIL"
$ws.Cells.Item(857,5).Value = "1. Video Streaming  (ms)"
$ws.Cells.Item(857,6).Value = "1. Static Yotube Streaming"
$ws.Cells.Item(857,7).Value = "1. Please update the AZQ app to version v3.2.822.apk specifically for the YouTube test.Kindly note that all other tests must continue to be performed using version v3.2.237.While performing the YouTube test, please ensure that the video is successfully playing in the script before saving the log file"
$fmtSrc.Copy()
$ws.Range("A857:G857").PasteSpecial(-4122)
$ws.Rows.Item(857).RowHeight = 104

$excel.CutCopyMode = $false

# Update selection to match the post-edit active cell
$ws.Range("E846").Select()

Write-Output "Added rows 846-857"
